$wb = $excel.ActiveWorkbook

# ---- Sheet: general ----
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 44.93688026761753
$ws.Range("B6").Value = 33.21688026761753
$ws.Range("B7").Value = 0.6217041096856285
$ws.Range("B8").Value = 0.6217041096856285
$ws.Range("B9").Value = 11.72
$ws.Range("B10").Value = 0

# ---- Sheet: alpha ----
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A2").Value = 9
$ws.Range("A3:A8").EntireRow.Delete() | Out-Null

# ---- Sheet: x ----
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 7
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 13
$ws.Range("B7").Value = 11
$ws.Range("B8").Value = 8
$ws.Range("B9").Value = 6
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 10
$ws.Range("B13").Value = 12
$ws.Range("B14").Value = 9

# ---- Sheet: U ----
$ws = $wb.Worksheets.Item("U")
$ws.Range("B3").Value = 3

# ---- Sheet: TBar ----
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 30.60033324079214
$ws.Range("B4").Value = 34.72107346555759
$ws.Range("B5").Value = 34.16886835983306
$ws.Range("B6").Value = 30
$ws.Range("B8").Value = 33.8200130889828
$ws.Range("B9").Value = 30
$ws.Range("B10").Value = 32.31224998648503
$ws.Range("B11").Value = 36.63062389744944
$ws.Range("B12").Value = 35.22705701268762
$ws.Range("B13").Value = 36.48759645946009
$ws.Range("B14").Value = 39.73914444587759
$ws.Range("B15").Value = 37.3015579161986

# ---- Sheet: y ----
$ws = $wb.Worksheets.Item("y")
$ws.Range("A2").Value = 9
$ws.Range("A3:A8").EntireRow.Delete() | Out-Null

# ---- Sheet: Q ----
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 193.0200000000017
$ws.Range("C8").Value = 202.3100000000017
$ws.Range("C9").Value = 191.2450000000017
$ws.Range("C10").Value = 208.9250000000017
$ws.Range("C11").Value = 197.6600000000017
$ws.Range("C12").Value = 274.4950000000024
$ws.Range("C13").Value = 282.9900000000024
$ws.Range("C14").Value = 275.9600000000024
$ws.Range("C15").Value = 289.3600000000025
$ws.Range("C16").Value = 285.0050000000024
$ws.Range("C17").Value = 224.1799999999995
$ws.Range("C18").Value = 224.6649999999995
$ws.Range("C19").Value = 201.1149999999995
$ws.Range("C20").Value = 218.9699999999995
$ws.Range("C21").Value = 207.1049999999995
$ws.Range("C22").Value = 72.6299999999995
$ws.Range("C23").Value = 80.0549999999995
$ws.Range("C24").Value = 82.31999999999948
$ws.Range("C25").Value = 83.9549999999995
$ws.Range("C26").Value = 80.8149999999995
$ws.Range("C32").Value = 177.26
$ws.Range("C33").Value = 188.5600000000001
$ws.Range("C34").Value = 166.2
$ws.Range("C35").Value = 181.48
$ws.Range("C36").Value = 167.6700000000001
$ws.Range("C37").Value = 83.07500000000087
$ws.Range("C38").Value = 84.72000000000087
$ws.Range("C39").Value = 77.97500000000086
$ws.Range("C40").Value = 89.77000000000088
$ws.Range("C41").Value = 81.85500000000087
$ws.Range("C42").Value = 140.5549999999989
$ws.Range("C43").Value = 159.2150000000007
$ws.Range("C44").Value = 142.1399999999989
$ws.Range("C45").Value = 147.7249999999989
$ws.Range("C46").Value = 139.7449999999989
$ws.Range("C47").Value = 283.9899999999987
$ws.Range("C48").Value = 305.86
$ws.Range("C49").Value = 283.5849999999987
$ws.Range("C50").Value = 299.1099999999986
$ws.Range("C51").Value = 277.3899999999987
$ws.Range("C52").Value = 198.5049999999982
$ws.Range("C53").Value = 217.895
$ws.Range("C54").Value = 203.8699999999982
$ws.Range("C55").Value = 208.3799999999981
$ws.Range("C56").Value = 192.6599999999982
$ws.Range("C57").Value = 224.1799999999995
$ws.Range("C58").Value = 224.6649999999995
$ws.Range("C59").Value = 201.1149999999995
$ws.Range("C60").Value = 218.9699999999995
$ws.Range("C61").Value = 207.1049999999995
$ws.Range("C62").Value = 283.9899999999987
$ws.Range("C63").Value = 305.86
$ws.Range("C64").Value = 283.5849999999987
$ws.Range("C65").Value = 299.1099999999986
$ws.Range("C66").Value = 277.3899999999987
$ws.Range("C67").Value = 274.4950000000024
$ws.Range("C68").Value = 282.9900000000024
$ws.Range("C69").Value = 275.9600000000024
$ws.Range("C70").Value = 289.3600000000025
$ws.Range("C71").Value = 285.0050000000024

# ---- Sheet: R ----
$ws = $wb.Worksheets.Item("R")
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C8").Value = 5.86
$ws.Range("C10").Value = 0

# ---- Sheet: L ----
$ws = $wb.Worksheets.Item("L")
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0

# ---- Sheet: rho ----
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2").Value = 9
$ws.Range("A3:A8").EntireRow.Delete() | Out-Null

Write-Host "edits applied"
